$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.6627635770477539
$ws.Range("AG1").Value = 0.8210322443515039
$ws.Range("D2").Value = 0.66011576012282414
$ws.Range("L2").Value = 0.98984623513661751
$ws.Range("A3").Value = 0.90843364261726056
$ws.Range("F4").Value = 0.61380730461921051
$ws.Range("C5").Value = 0.57343258098852745
$ws.Range("D5").Value = 0.99545879528352765
$ws.Range("E6").Value = 0.92056352618067283
$ws.Range("AR6").Value = 0.91065607495144119
$ws.Range("E7").Value = 0.69516830662048679
$ws.Range("F7").Value = 0.58033161882477435
$ws.Range("AP7").Value = 0.97312102141290135
$ws.Range("F8").Value = 0.95164849122104833
$ws.Range("G8").Value = 0.83804205525964703
$ws.Range("I8").Value = 0.9205081835231399
$ws.Range("J8").Value = 0.5537710738785544
$ws.Range("I10").Value = 0.65146918175096702
$ws.Range("K10").Value = 0.81059055348648357
$ws.Range("I11").Value = 0.89177981254379191
$ws.Range("K12").Value = 0.93759434874326075
$ws.Range("BH12").Value = 0.71643963516056031
$ws.Range("K13").Value = 0.94332644585749392
$ws.Range("N13").Value = 0.60658721813827743
$ws.Range("O13").Value = 0.98678034303218654
$ws.Range("L14").Value = 0.70866289039196839
$ws.Range("N15").Value = 0.82670088066492775
$ws.Range("Q15").Value = 0.63223241838303523
$ws.Range("C16").Value = 0.71076852777530086
$ws.Range("N16").Value = 0.92568661111378447
$ws.Range("Q16").Value = 0.97568559589718262
$ws.Range("T18").Value = 0.61308306644825228
$ws.Range("Q19").Value = 0.85890497409347932
$ws.Range("T19").Value = 0.88548754781490557
$ws.Range("X19").Value = 0.87894017791781676
$ws.Range("V20").Value = 0.8667875398720879
$ws.Range("G21").Value = 0.9737342824088302
$ws.Range("T21").Value = 0.97417417515951388
$ws.Range("U22").Value = 0.90173862831076934
$ws.Range("W22").Value = 0.64468279773040094
$ws.Range("AM22").Value = 0.95898233038739455
$ws.Range("AA23").Value = 0.61962277765031626
$ws.Range("AN23").Value = 0.99335128143495255
$ws.Range("T24").Value = 0.72430745344395553
$ws.Range("Z25").Value = 0.88032390248100456
$ws.Range("AA25").Value = 0.55455651533137518
$ws.Range("Z27").Value = 0.82931245072500692
$ws.Range("AC27").Value = 0.97279598030263648
$ws.Range("Z28").Value = 0.95189235385746174
$ws.Range("AC28").Value = 0.89154211017254603
$ws.Range("AG29").Value = 0.94755591763839386
$ws.Range("AA30").Value = 0.83505433468796031
$ws.Range("AB30").Value = 0.64472089055923465
$ws.Range("AE30").Value = 0.89399839375890633
$ws.Range("O31").Value = 0.98645496465739568
$ws.Range("AC31").Value = 0.82063226762899544
$ws.Range("AF31").Value = 0.93024928710132604
$ws.Range("AE33").Value = 0.88819550585779061
$ws.Range("AF33").Value = 0.96040692518841952
$ws.Range("AF34").Value = 0.74042129249552624
$ws.Range("AG35").Value = 0.92927719302643363
$ws.Range("AH35").Value = 0.67667885461628396
$ws.Range("AK35").Value = 0.97445493005782902
$ws.Range("AH36").Value = 0.55406528049790116
$ws.Range("AI36").Value = 0.99937002814028508
$ws.Range("AK36").Value = 0.81485964913482001
$ws.Range("AL37").Value = 0.86741640386458474
$ws.Range("V38").Value = 0.81503683109307401
$ws.Range("W38").Value = 0.87014975569356079
$ws.Range("AJ38").Value = 0.86106657065543291
$ws.Range("AM38").Value = 0.88770673595576366
$ws.Range("BL38").Value = 0.8863071566648093
$ws.Range("Q39").Value = 0.98565361625036907
$ws.Range("R39").Value = 0.90189724457502018
$ws.Range("AK39").Value = 0.88353770782315055
$ws.Range("AM40").Value = 0.95513611172410617
$ws.Range("AO40").Value = 0.98850036577421729
$ws.Range("AP40").Value = 0.53619216225254274
$ws.Range("O41").Value = 0.98524205889452654
$ws.Range("AQ41").Value = 0.96096494922518305
$ws.Range("BA42").Value = 0.98090648109362322
$ws.Range("AP43").Value = 0.98476763954830515
$ws.Range("AR43").Value = 0.99466568502400154
$ws.Range("AS43").Value = 0.81015072691521084
$ws.Range("AP44").Value = 0.61360725863380128
$ws.Range("AR45").Value = 0.88455090224649546
$ws.Range("AS46").Value = 0.97012270953835777
$ws.Range("AU46").Value = 0.76579263635030581
$ws.Range("AV46").Value = 0.67061571407939136
$ws.Range("AS47").Value = 0.98573509228970202
$ws.Range("AW47").Value = 0.89556673126028197
$ws.Range("O48").Value = 0.57273828954621031
$ws.Range("AU48").Value = 0.87481315759691136
$ws.Range("AV49").Value = 0.9418393013092774
$ws.Range("AX49").Value = 0.89415054010109141
$ws.Range("C50").Value = 0.90529010713498559
$ws.Range("AW51").Value = 0.89733465724973249
$ws.Range("AX51").Value = 0.92307723923568019
$ws.Range("BC51").Value = 0.81316148557670598
$ws.Range("Y52").Value = 0.89620592346942907
$ws.Range("AY53").Value = 0.68134012238965669
$ws.Range("AZ54").Value = 0.60723441897053965
$ws.Range("BH54").Value = 0.89966093833835381
$ws.Range("AH55").Value = 0.67283204474616687
$ws.Range("BB55").Value = 0.71376106258206751
$ws.Range("B56").Value = 0.97786579686798158
$ws.Range("BF56").Value = 0.84714290050449415
$ws.Range("BC57").Value = 0.91328172857647116
$ws.Range("BD57").Value = 0.91591248180193419
$ws.Range("BG57").Value = 0.91235482759997311
$ws.Range("AG58").Value = 0.98198667144841689
$ws.Range("BE58").Value = 0.61422810043633813
$ws.Range("BG58").Value = 0.87687843500088958
$ws.Range("BD59").Value = 0.8703713827009274
$ws.Range("BH59").Value = 0.64014314537040562
$ws.Range("S61").Value = 0.55794364463306834
$ws.Range("BH61").Value = 0.7303681977199743
$ws.Range("BK61").Value = 0.65870468115148628
$ws.Range("I62").Value = 0.70231409207301754
$ws.Range("BK62").Value = 0.85600349492646943
$ws.Range("X63").Value = 0.96689448985332049
$ws.Range("BL63").Value = 0.83994271303276768
$ws.Range("BN64").Value = 0.88014190828966588
$ws.Range("V65").Value = 0.8258685659315641
$ws.Range("BK65").Value = 0.65260239171286938
$ws.Range("BM66").Value = 0.98817290592966611
$ws.Range("BO66").Value = 0.95093582540825283
$ws.Range("BP66").Value = 0.93085686534383705
$ws.Range("A67").Value = 0.80020593135401552
$ws.Range("J67").Value = 0.7076862800795487
$ws.Range("P67").Value = 0.96553591759299151
$ws.Range("BA67").Value = 0.91978417302280202
$ws.Range("BM67").Value = 0.90171775904421936
$ws.Range("B68").Value = 0.71867225709792604
$ws.Range("AU68").Value = 0.99287599909831981
$ws.Range("BO68").Value = 0.80583498786211138
